# Refresh the crypto price/volume snapshot (GitHub Actions scheduled update).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Column D ("Price") holds plain text in the source sheet (e.g. "211.10", "1.711.49")
# -- several of these look numeric to Excel and would otherwise be auto-coerced to a
# float (dropping trailing zeros / punctuation). Prefix new Price values with a
# leading apostrophe, Excel's normal "force text" quote-prefix marker, so they are
# stored verbatim as text, matching the original formatting exactly.
function Set-TextValue($range, [string]$text) {
    $ws.Range($range).Value = "'" + $text
}

Set-TextValue "D2" '28.403.99'
$ws.Range("E2").Value = '  +1.90%  '
Set-TextValue "D3" '1.573.73'
$ws.Range("E3").Value = '  +0.22%  '
$ws.Range("E4").Value = '  +1.30%  '
Set-TextValue "D5" '211.10'
$ws.Range("E5").Value = '  -0.03%  '
$ws.Range("E6").Value = '  -0.93%  '
$ws.Range("E7").Value = '  +1.13%  '
Set-TextValue "D8" '45.97'
$ws.Range("E8").Value = '  +4.01%  '
Set-TextValue "D9" '23.71'
$ws.Range("E9").Value = '  +2.29%  '
$ws.Range("E10").Value = '  -1.11%  '
Set-TextValue "D11" '0.0591'
$ws.Range("E11").Value = '  -0.98%  '
$ws.Range("E12").Value = '  +0.44%  '
Set-TextValue "D13" '1.799.50'
$ws.Range("E13").Value = '  +0.32%  '
Set-TextValue "D14" '1.574.58'
$ws.Range("E14").Value = '  +0.26%  '
Set-TextValue "D15" '0.522'
$ws.Range("E15").Value = '  +0.42%  '
$ws.Range("E16").Value = '  -1.22%  '
Set-TextValue "D17" '28.396.57'
$ws.Range("E17").Value = '  +2.13%  '
Set-TextValue "D18" '62.26'
$ws.Range("E18").Value = '  -1.91%  '
Set-TextValue "D19" '227.92'
$ws.Range("E19").Value = '  -1.16%  '
Set-TextValue "D20" '7.35'
$ws.Range("E20").Value = '  -1.59%  '
Set-TextValue "D21" '0.0₃0692'
$ws.Range("E21").Value = '  -1.47%  '
$ws.Range("E22").Value = '  +1.12%  '
Set-TextValue "D23" '3.92'
$ws.Range("E23").Value = '  -4.54%  '
Set-TextValue "D24" '9.16'
$ws.Range("E24").Value = '  -1.55%  '
$ws.Range("E25").Value = '  +3.26%  '
Set-TextValue "D26" '150.61'
$ws.Range("E26").Value = '  +0.10%  '
Set-TextValue "D27" '14.98'
$ws.Range("E27").Value = '  -1.71%  '
Set-TextValue "D28" '6.45'
$ws.Range("E28").Value = '  -1.72%  '
$ws.Range("E29").Value = '  -2.49%  '
$ws.Range("E30").Value = '  +1.16%  '
Set-TextValue "D31" '1.11'
$ws.Range("E31").Value = '  -1.92%  '
$ws.Range("E32").Value = '  -1.82%  '
$ws.Range("E33").Value = '  -0.98%  '
Set-TextValue "D34" '3.11'
$ws.Range("E34").Value = '  -0.55%  '
Set-TextValue "D35" '1.393.13'
$ws.Range("E35").Value = '  -1.67%  '
$ws.Range("E36").Value = '  -2.35%  '
Set-TextValue "D37" '0.998'
$ws.Range("E37").Value = '  -3.75%  '
$ws.Range("E38").Value = '  +3.03%  '
Set-TextValue "D39" '2.55'
$ws.Range("E39").Value = '  +3.71%  '
$ws.Range("E40").Value = '  -1.07%  '
Set-TextValue "D41" '0.531'
$ws.Range("E41").Value = '  -2.05%  '
$ws.Range("E42").Value = '  +1.18%  '
$ws.Range("E43").Value = '  -1.67%  '
Set-TextValue "D44" '5.62'
$ws.Range("E44").Value = '  -0.68%  '
$ws.Range("E45").Value = '  +0.69%  '
$ws.Range("E46").Value = '  +1.60%  '
Set-TextValue "D47" '62.20'
$ws.Range("E47").Value = '  -2.45%  '
Set-TextValue "D48" '1.711.63'
$ws.Range("E48").Value = '  +0.83%  '
Set-TextValue "D49" '85.65'
$ws.Range("E49").Value = '  -1.14%  '
$ws.Range("B50").Value = 'BabyDogeCoin'
$ws.Range("C50").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
Set-TextValue "D50" '0.0₆0101'
$ws.Range("E50").Value = '  +2.09%  '
$ws.Range("B51").Value = 'Cronos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextValue "D51" '0.0519'
$ws.Range("E51").Value = '  -0.72%  '
